$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column D: "expectedwear" header + default values ---
$ws.Range("D1").Value = "expectedwear"
$ws.Range("D2:D6").Value = 14
$ws.Range("D7:D13").Value = 10

# --- Column A width change (21.5 -> ~32.33 chars) ---
$ws.Columns.Item(1).ColumnWidth = 31.5

# --- New row 14: styled (bold header-like) empty cell A14 ---
$c = $ws.Range("A14")
$c.NumberFormat = "@"
$c.Font.Name = "Helvetica Neue"
$c.Font.Size = 10
$c.Font.Bold = $true
$c.Font.ColorIndex = 1
$c.Interior.Pattern = 1
$c.Interior.ColorIndex = 2
$c.Borders.LineStyle = 1
$c.Borders.Item(7).ColorIndex = 3
$c.Borders.Item(10).ColorIndex = 3
$c.Borders.Item(8).ColorIndex = 3
$c.Borders.Item(9).ColorIndex = 4
$c.VerticalAlignment = -4160

# --- Selection / view state ---
$ws.Range("D15").Select()
